$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so they stay text (matches source inlineStr type)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = "44.105.48"
$ws.Range("E2").Value = "  +4.47%  "
$ws.Range("D3").Value = "2.218.35"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "259.76"
$ws.Range("E5").Value = "  +2.64%  "
$ws.Range("D6").Value = "82.35"
$ws.Range("E6").Value = "  +12.22%  "
$ws.Range("E7").Value = "  +3.20%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.608"
$ws.Range("E9").Value = "  +5.26%  "
$ws.Range("D10").Value = "44.09"
$ws.Range("E10").Value = "  +8.63%  "
$ws.Range("D11").Value = "0.0932"
$ws.Range("E11").Value = "  +2.53%  "
$ws.Range("D12").Value = "7.05"
$ws.Range("E12").Value = "  +3.93%  "
$ws.Range("E13").Value = "  +3.14%  "
$ws.Range("D14").Value = "2.556.29"
$ws.Range("E14").Value = "  +2.26%  "
$ws.Range("D15").Value = "14.61"
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").Value = "2.230.73"
$ws.Range("E16").Value = "  +3.01%  "
$ws.Range("D17").Value = "0.782"
$ws.Range("E17").Value = "  +2.89%  "
$ws.Range("D18").Value = "43.972.74"
$ws.Range("E18").Value = "  +4.43%  "
$ws.Range("E19").Value = "  +1.68%  "
$ws.Range("D20").Value = "71.46"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").Value = "6.01"
$ws.Range("E21").Value = "  +2.75%  "
$ws.Range("E22").Value = "  +9.17%  "
$ws.Range("D23").Value = "233.05"
$ws.Range("E23").Value = "  +3.37%  "
$ws.Range("D24").Value = "9.30"
$ws.Range("E24").Value = "  -2.23%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "10.78"
$ws.Range("E26").Value = "  +3.09%  "
$ws.Range("D27").Value = "41.34"
$ws.Range("E27").Value = "  +12.64%  "
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("E29").Value = "  +2.71%  "
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").Value = "  +3.10%  "
$ws.Range("D31").Value = "172.85"
$ws.Range("E31").Value = "  +2.32%  "
$ws.Range("D32").Value = "20.65"
$ws.Range("E32").Value = "  +3.21%  "
$ws.Range("D33").Value = "0.0883"
$ws.Range("E33").Value = "  +10.08%  "
$ws.Range("D34").Value = "5.33"
$ws.Range("E34").Value = "  +4.15%  "
$ws.Range("D35").Value = "0.116"
$ws.Range("E35").Value = "  +7.77%  "
$ws.Range("E36").Value = "  +1.99%  "
$ws.Range("D37").Value = "0.0362"
$ws.Range("E37").Value = "  +9.62%  "
$ws.Range("D38").Value = "4.49"
$ws.Range("E38").Value = "  +6.18%  "
$ws.Range("D39").Value = "13.41"
$ws.Range("E39").Value = "  +12.11%  "
$ws.Range("D40").Value = "2.97"
$ws.Range("E40").Value = "  +21.09%  "
$ws.Range("E41").Value = "  +3.44%  "
$ws.Range("D42").Value = "63.50"
$ws.Range("E42").Value = "  +7.69%  "
$ws.Range("D43").Value = "5.55"
$ws.Range("E43").Value = "  +8.85%  "
$ws.Range("D44").Value = "0.201"
$ws.Range("E44").Value = "  +3.00%  "
$ws.Range("D45").Value = "102.86"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("D46").Value = "0.0990"
$ws.Range("E46").Value = "  +2.21%  "
$ws.Range("D47").Value = "8.33"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "1.12"
$ws.Range("E48").Value = "  +3.42%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "1.56"
$ws.Range("E49").Value = "  +27.79%  "
$ws.Range("E50").Value = "  +3.64%  "
$ws.Range("D51").Value = "0.442"
$ws.Range("E51").Value = "  -5.19%  "
